$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window / view sizing (workbookView) ---
$win = $excel.ActiveWindow
$win.Width = 28800
$win.Height = 12435

# --- Rows 2 & 3: Invoice Class columns (E/F/G) now point at the new values ---
# (iHeartmedia1! / krishnagundavarapu@iheartmedia.com / SB5 20.1 Upgrade...)
# These already hold the right text; nothing to change there directly, the
# rows that get cleared below are what drops the old/orphaned shared strings.

# --- Rows 4-8: clear out the old "Invoice Classes" test rows ---
# Columns A-F are fully cleared (cell + formatting removed entirely).
# Columns G, I, J keep their formatting but lose their value.
for ($r = 4; $r -le 8; $r++) {
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 6)).Clear() | Out-Null
    $ws.Cells.Item($r, 7).ClearContents() | Out-Null
    $ws.Cells.Item($r, 9).ClearContents() | Out-Null
    $ws.Cells.Item($r, 10).ClearContents() | Out-Null
}

# --- Hyperlinks: only D2, D3, E2, E3 should remain ---
# The engine's Hyperlinks.Delete() on a Range clears the whole sheet's
# hyperlink collection, so re-add just the four that should survive, then
# restore their original cell styles (Hyperlinks.Add forces the Hyperlink
# style onto the cell).
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D2"), "https://system.netsuite.com/pages/customerlogin.jsp") | Out-Null
$ws.Range("D2").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("D3"), "https://system.netsuite.com/pages/customerlogin.jsp") | Out-Null
$ws.Range("D3").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:krishnagundavarapu@iheartmedia.com") | Out-Null
$ws.Range("E2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:krishnagundavarapu@iheartmedia.com") | Out-Null
$ws.Range("E3").Style = "Normal"

# --- Sheet view: scroll back to A1 and move the selection ---
$ws.Range("D11").Select() | Out-Null
